# ducvh cap nhat hdct, dm, sp + writeexcelv3
#
# New batch of "fail date" product rows (ID 50-54, newer MaSP/SoLuong/
# HanSuDung values) is inserted at the top of Sheet1, ahead of the
# previously-recorded rows. The very first legacy row (the original
# ID=50 / MaSP=SP01 record) is superseded/removed, while the other
# legacy rows (originally rows 2-5) shift down beneath the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the 5 existing rows down by 5 so there is room for the new batch
# (old row1..row5 -> row6..row10).
$ws.Rows("1:5").Insert()

# Write the new product rows into the freshly opened rows 1-5.
$ws.Range("A1").Value = "ID= 50;IDLS= 150;MaSP= SP041;SoLuong= 978;HanSuDung= 18-02-2022;TrangThai= true"
$ws.Range("A2").Value = "ID= 51;IDLS= 151;MaSP= SP027;SoLuong= 346;HanSuDung= 21-04-2022;TrangThai= true"
$ws.Range("A3").Value = "ID= 52;IDLS= 152;MaSP= SP023;SoLuong= 866;HanSuDung= 19-10-2022;TrangThai= true"
$ws.Range("A4").Value = "ID= 53;IDLS= 153;MaSP= SP015;SoLuong= 383;HanSuDung= 13-09-2022;TrangThai= true"
$ws.Range("A5").Value = "ID= 54;IDLS= 154;MaSP= SP024;SoLuong= 970;HanSuDung= 13-10-2022;TrangThai= true"

# The old first record (now pushed to row 6 - "ID=50;...;MaSP=SP01;...")
# is dropped entirely; the remaining legacy rows close up beneath the
# new rows (old row2..row5 end up as row6..row9).
$ws.Rows(6).Delete()

$wb.Save()
